# Handback report generation:
#  - Status text updated from "Ready for handoff" to "Handed back: in sync with en-US"
#    everywhere it appears (Overview + per-locale sheets).
#  - "Latest Target File" (col I) and "Latest Handback File" (col J) populated for
#    both locale sheets (zh-cn, de-de), with I2/I3 turned into hyperlinks (same
#    look as the existing A2/A3/B2/B3 hyperlinks).
#  - "Latest Handback DateTime" (col K) timestamps filled in.
#  - Columns I/J/C widened to fit the new content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# blue hyperlink font color (matches the workbook's existing custom "HyperLink" style)
$hyperlinkColor = 15570276

# ---- Status text (shared across sheets) ----
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---- zh-cn: target/handback file + datetime ----
$zhTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6f540da7c0ad0eb0d18922bc8e97749f485d043/e2e/a.md"
$zhSourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6f540da7c0ad0eb0d18922bc8e97749f485d043/e2e/b.md"

$wsZhCn.Range("I2").Value = "a.md"
$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-02 10:43:59"

$wsZhCn.Range("I3").Value = "a.md"
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-02 10:43:59"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhTargetUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhTargetUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhSourceUrl, [Type]::Missing, [Type]::Missing, "b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhTargetUrl, [Type]::Missing, [Type]::Missing, "a.md")

foreach ($addr in @("A2", "I2", "A3", "I3")) {
    $r = $wsZhCn.Range($addr)
    $r.Font.Name = "Calibri"
    $r.Font.Underline = 2
    $r.Font.Color = $hyperlinkColor
}

$wsZhCn.Columns("C").ColumnWidth = 29.144371396019366
$wsZhCn.Columns("J").ColumnWidth = 39.166666666666664

# ---- de-de: target/handback file + datetime ----
$deTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6f540da7c0ad0eb0d18922bc8e97749f485d043/e2e/a.md"
$deSourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6f540da7c0ad0eb0d18922bc8e97749f485d043/e2e/b.md"

$wsDeDe.Range("I2").Value = "a.md"
$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-02 10:44:14"

$wsDeDe.Range("I3").Value = "a.md"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-02 10:44:14"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deTargetUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $deTargetUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $deSourceUrl, [Type]::Missing, [Type]::Missing, "b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $deTargetUrl, [Type]::Missing, [Type]::Missing, "a.md")

foreach ($addr in @("A2", "I2", "A3", "I3")) {
    $r = $wsDeDe.Range($addr)
    $r.Font.Name = "Calibri"
    $r.Font.Underline = 2
    $r.Font.Color = $hyperlinkColor
}

$wsDeDe.Columns("C").ColumnWidth = 29.144371396019366
$wsDeDe.Columns("J").ColumnWidth = 39.166666666666664

# ---- Overview: widen Status columns ----
$wsOverview.Columns("E").ColumnWidth = 29.144371396019366
$wsOverview.Columns("F").ColumnWidth = 29.144371396019366
